$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date column (A2:A7) from 2025-12-05 to 2025-12-06.
# Force text formatting first so Excel doesn't auto-convert the
# "yyyy-mm-dd" looking text into a real date serial number, then
# restore the default "Normal" style so no extra formatting sticks
# to the cells.
$ws.Range("A2:A7").NumberFormat = "@"
$ws.Range("A2:A7").Value = "2025-12-06"
$ws.Range("A2:A7").Style = "Normal"

# Update 최종점수 (column K, "final score") values
$ws.Range("K2").Value = 62.7
$ws.Range("K3").Value = 54.5
$ws.Range("K4").Value = 50.7
$ws.Range("K5").Value = 47.9
$ws.Range("K6").Value = 39.9
$ws.Range("K7").Value = 39.9

# Update MACRO_SCORE (column N) values for all data rows
$ws.Range("N2:N7").Value = 51.54219175917372
